$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Remove the narrow spacer columns (D, H, L) that separated the three
# per-cluster metric blocks, shifting the remaining data left.
$ws.Range("L:L").Delete()
$ws.Range("H:H").Delete()
$ws.Range("D:D").Delete()

# Update the active selection to match the saved state after the edit.
$ws.Range("Q12").Select()

$wb.RunCommand("review.inquire.clean", @{ sheet = $ws.Name })
